# Update the SQL JOIN clauses in the StatQuery/TabQuery cells (B2:B7, C2) on Sheet1.
# The author renamed the join keys from the generic "id" columns to the
# fully-qualified "<table>_id" columns (e.g. std.id -> std.study_id,
# prt.id -> prt.participant_id) across every embedded query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of (old substring, new substring) applied, in order, to every query cell.
$replacements = @(
    @('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"'),
    @('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'),
    @('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'),
    @('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'),
    @('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"'),
    @('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
)

$targetCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellAddr in $targetCells) {
    $cell = $ws.Range($cellAddr)
    $text = $cell.Value2
    foreach ($pair in $replacements) {
        $old = $pair[0]
        $new = $pair[1]
        # Escape the search text as a regex literal; double any `$` in the
        # replacement so .NET regex doesn't treat it as a backreference token.
        $safeNew = $new.Replace('$', '$$')
        $text = $text -replace [regex]::Escape($old), $safeNew
    }
    $cell.Value = $text
}
